$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize oxidation-state strings so the negative value is listed first.
# Several rows share the same underlying text, so every row that displayed
# the old value must be updated to the new one.

# Arsenico (As) / Nitrogeno (N) / Fosforo (P): "3,-3,5" -> "-3,3,5"
$ws.Range("D4").Value = "-3,3,5"
$ws.Range("D29").Value = "-3,3,5"
$ws.Range("D33").Value = "-3,3,5"

# Bromo (Br): "1,3,5,-1" -> "-1,1,3,5"
$ws.Range("D9").Value = "-1,1,3,5"

# Carbono (C): "2,4,-4" -> "-4,2,4"
$ws.Range("D10").Value = "-4,2,4"

# Cloro (Cl) / Yodo (I): "1,3,5,7-1" -> "-1,1,3,5,7"
$ws.Range("D13").Value = "-1,1,3,5,7"
$ws.Range("D24").Value = "-1,1,3,5,7"

# Hidrogeno (H): "1,-1" -> "-1,1"
$ws.Range("D22").Value = "-1,1"

# Azufre (S): "2,4,6,-2" -> "-2,2,4,6"
$ws.Range("D39").Value = "-2,2,4,6"

# Selenio (Se): corrected oxidation states "2,4,6,-2" -> "-2,2,4,7"
$ws.Range("D40").Value = "-2,2,4,7"

# Silicio (Si): "4,-4" -> "-4,4"
$ws.Range("D41").Value = "-4,4"

# Reflect where the editor finished up: selection moved to D47.
$ws.Range("D47").Select()
